$d = $word.ActiveDocument

# The last paragraph in the document body (before the sectPr) is the one
# containing "Quer dizer: O significado dos vocábulos...". In the edit, its
# paragraph-mark run properties (w:pPr/w:rPr) lose <w:b/> and <w:bCs/>, and
# six new paragraphs are appended after it (still before the final sectPr).
# We replace that paragraph's Range contents in one shot via InsertXML with
# the exact target OOXML for all seven paragraphs (the corrected "Quer
# dizer..." paragraph -- with its original w14:paraId/rsid attributes kept
# intact -- plus the six new ones). This both strips the bold from the
# paragraph mark and appends the new content in the right place.

$lastPara = $d.Paragraphs.Last
$targetRange = $lastPara.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="42EFBDA9" w14:textId="5144CC93" w:rsidR="00C434F6" w:rsidRPr="00C434F6" w:rsidRDefault="00C434F6" w:rsidP="00C434F6">
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:ind w:left="1080"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">i. </w:t>
  </w:r>
  <w:r w:rsidR="00F76D84">
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Quer dizer: O significado dos vocábulos, por oposição a sua forma, ou seja, tem muito mais valor o significado do que a forma.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:tab/>
    <w:t xml:space="preserve">***Podemos colocar uma parte de um texto entre tags automaticamente selecionando o texto e apertando </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>CTRL+SHIFT+P</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> e escolhendo a opção wrap “with abreviation” e digitando a tag que deseja ser colocada***</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Formatação de Textos:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Corpodetexto"/>
  </w:pPr>
  <w:r>
    <w:t xml:space="preserve">Existem diversos métodos de </w:t>
  </w:r>
  <w:r>
    <w:t>se formatar textos, o primeiro é colocar tags de formatação apenas, já o segundo, utiliza-se de tags semânticas para formatar os textos. Exemplos</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> de tags sem e com semântica, respectivamente</w:t>
  </w:r>
  <w:r>
    <w:t>:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="14"/>
    </w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>&lt;b&gt;/&lt;strong&gt;:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> Ambos deixam o texto em negrito.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="14"/>
    </w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>&lt;i&gt;/&lt;em&gt;:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> Ambos formatam o texto em itálico.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="PargrafodaLista"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="14"/>
    </w:numPr>
    <w:spacing w:line="360" w:lineRule="auto"/>
    <w:jc w:val="both"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve">&lt;u&gt;/: </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Sublinha o texto.</w:t>
  </w:r>
</w:p>
'@

$targetRange.InsertXML($xml)
